$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 75.3940389498535
$ws.Cells.Item(2,4).Value = -66613.2146299122
$ws.Cells.Item(2,5).Value = 133381.131621427
$ws.Cells.Item(2,7).Value = 1
$ws.Cells.Item(2,8).Value = 133977.992975931
$ws.Cells.Item(2,9).Value = 401.054344111106
$ws.Cells.Item(2,10).Value = 16509.6059610501

$ws.Cells.Item(3,1).Value = "aerobic_scope ~ s(doy, by = fish_basin, bs = `"cc`", k = 17) + s(floy_tag, year, by = fish_basin, bs = c(`"re`", `"re`"), k = c(20, 4)) + ACF"
$ws.Cells.Item(3,2).Value = "m21"
$ws.Cells.Item(3,3).Value = 62.7084781667592
$ws.Cells.Item(3,4).Value = -66762.0660903263
$ws.Cells.Item(3,5).Value = 133652.438076203
$ws.Cells.Item(3,6).Value = 271.306454775709
$ws.Cells.Item(3,7).Value = [double]"1.22053965217866e-59"
$ws.Cells.Item(3,8).Value = 134147.458512967
$ws.Cells.Item(3,9).Value = 408.780894170596
$ws.Cells.Item(3,10).Value = 16522.2915218332

$ws.Cells.Item(4,3).Value = 62.7241569845925
$ws.Cells.Item(4,4).Value = -66762.1072299664
$ws.Cells.Item(4,5).Value = 133652.548544725
$ws.Cells.Item(4,6).Value = 271.416923298122
$ws.Cells.Item(4,7).Value = [double]"1.15495206111892e-59"
$ws.Cells.Item(4,8).Value = 134147.677739164
$ws.Cells.Item(4,9).Value = 408.78049079072
$ws.Cells.Item(4,10).Value = 16522.2758430154

$ws.Cells.Item(5,1).Value = "aerobic_scope ~ fish_basin + s(doy, by = fish_basin, bs = `"cc`", k = 17) + s(floy_tag, year, by = fish_basin, bs = c(`"re`", `"re`"), k = c(20, 4)) + ti(doy, fish_basin, bs = c(`"cc`", `"fs`"), k = c(20, 3))"
$ws.Cells.Item(5,2).Value = "m8"
$ws.Cells.Item(5,3).Value = 77.4466570389077
$ws.Cells.Item(5,4).Value = -67007.1455135324
$ws.Cells.Item(5,5).Value = 134172.792312448
$ws.Cells.Item(5,6).Value = 791.660691021068
$ws.Cells.Item(5,7).Value = [double]"1.238982507667e-172"
$ws.Cells.Item(5,8).Value = 134784.310397268
$ws.Cells.Item(5,9).Value = 395.702420303694
$ws.Cells.Item(5,10).Value = 16507.5533429611

$ws.Cells.Item(6,6).Value = 1121.26266361601
$ws.Cells.Item(6,7).Value = [double]"3.31822792381082e-244"

$ws.Cells.Item(7,6).Value = 1121.27805216474
$ws.Cells.Item(7,7).Value = [double]"3.29279453896598e-244"

$ws.Cells.Item(8,6).Value = 2442.83075221756

$ws.Cells.Item(9,6).Value = 8572.41379386108

$ws.Cells.Item(10,6).Value = 8587.46555954442

$ws.Cells.Item(11,6).Value = 8587.86125133466

$ws.Cells.Item(12,6).Value = 8656.7392633337

$ws.Cells.Item(13,6).Value = 8739.8991834491

$ws.Cells.Item(14,6).Value = 8755.90258363585

$ws.Cells.Item(15,6).Value = 9503.03232447105

$ws.Cells.Item(16,6).Value = 9503.03232447105

$ws.Cells.Item(17,6).Value = 20557.0785581974

$ws.Cells.Item(18,6).Value = 24019.2404132324
